$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Seguimiento")

# Fill in the dates for row 6 (C6:E6) matching B6's existing date style
# Serial date 42068 = March 5, 2015
$ws.Range("C6").Value = 42068
$ws.Range("D6").Value = 42068
$ws.Range("E6").Value = 42068

# Update the active selection from C6 to E6
$ws.Range("E6").Select()
